$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT, preserving the cell's original style
# (Excel's COM Value setter auto-coerces plain-looking numeric strings,
# e.g. "213.25" or "23.06", into floating point numbers. Forcing the
# NumberFormat to "@" (text) before the assignment keeps it a string,
# then restoring the original Style avoids leaving any visible
# formatting change behind.)
function Set-TextValue($range, [string]$value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "27.702.46"
$ws.Range("E2").Value = "  +1.30%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "1.646.58"
$ws.Range("E3").Value = "  -0.58%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.07%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "213.25"

# Row 6 - XRP
Set-TextValue $ws.Range("D6") "0.531"
$ws.Range("E6").Value = "  +3.03%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.06%  "

# Row 8 - Solana
Set-TextValue $ws.Range("D8") "23.06"
$ws.Range("E8").Value = "  -2.63%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -1.35%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.35%  "

# Row 11 - TRON
Set-TextValue $ws.Range("D11") "0.0888"
$ws.Range("E11").Value = "  +1.43%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D12") "1.879.13"
$ws.Range("E12").Value = "  -0.60%  "

# Row 13 - WrappedEther
Set-TextValue $ws.Range("D13") "1.641.48"
$ws.Range("E13").Value = "  -0.86%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -0.93%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  -1.60%  "

# Row 16 - Litecoin
Set-TextValue $ws.Range("D16") "64.20"
$ws.Range("E16").Value = "  -2.50%  "

# Row 17 - WrappedBTC
Set-TextValue $ws.Range("D17") "27.664.86"
$ws.Range("E17").Value = "  +1.14%  "

# Row 18 - BitcoinCash
Set-TextValue $ws.Range("D18") "230.14"
$ws.Range("E18").Value = "  -0.79%  "

# Row 19 - ShibaInu
Set-TextValue $ws.Range("D19") "0.0₃0725"
$ws.Range("E19").Value = "  -0.31%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  +2.00%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.08%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -1.38%  "

# Row 23 - Avalanche
Set-TextValue $ws.Range("D23") "10.02"
$ws.Range("E23").Value = "  +6.72%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -3.26%  "

# Row 25 - Monero
Set-TextValue $ws.Range("D25") "149.01"
$ws.Range("E25").Value = "  +1.29%  "

# Row 26 - Cosmos
$ws.Range("E26").Value = "  -2.79%  "

# Row 27 - Stellar
Set-TextValue $ws.Range("D27") "0.112"
$ws.Range("E27").Value = "  +0.67%  "

# Row 28 - BinanceUSD
$ws.Range("E28").Value = "  -0.01%  "

# Row 29 - EthereumClassic
Set-TextValue $ws.Range("D29") "15.65"

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -0.02%  "

# Row 31 - Hedera
Set-TextValue $ws.Range("D31") "0.0484"
$ws.Range("E31").Value = "  -2.66%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -0.06%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D33") "3.18"
$ws.Range("E33").Value = "  +1.75%  "

# Row 34 - Maker
Set-TextValue $ws.Range("D34") "1.441.68"
$ws.Range("E34").Value = "  -0.72%  "

# Row 35 - LidoDAOToken
Set-TextValue $ws.Range("D35") "1.59"
$ws.Range("E35").Value = "  +1.46%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  -1.16%  "

# Row 37 - ImmutableX
$ws.Range("E37").Value = "  -0.05%  "

# Row 38 - ARBITRUM
Set-TextValue $ws.Range("D38") "0.883"
$ws.Range("E38").Value = "  -2.92%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  -1.11%  "

# Row 40 - TrustWalletToken
Set-TextValue $ws.Range("D40") "0.903"
$ws.Range("E40").Value = "  +15.01%  "

# Row 41 - WEMIXToken
$ws.Range("E41").Value = "  -0.92%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  +0.13%  "

# Row 43 - FraxShare
Set-TextValue $ws.Range("D43") "5.69"
$ws.Range("E43").Value = "  +3.11%  "

# Row 44 - mCoin
$ws.Range("E44").Value = "  -0.06%  "

# Row 45 - MXToken
$ws.Range("E45").Value = "  +1.80%  "

# Row 46 - Aave
Set-TextValue $ws.Range("D46") "65.61"
$ws.Range("E46").Value = "  +0.51%  "

# Row 47 - RocketPoolETH
Set-TextValue $ws.Range("D47") "1.788.60"
$ws.Range("E47").Value = "  -0.48%  "

# Row 48 - RenderToken
Set-TextValue $ws.Range("D48") "1.69"
$ws.Range("E48").Value = "  -1.47%  "

# Row 49 - Quant
Set-TextValue $ws.Range("D49") "86.41"
$ws.Range("E49").Value = "  -2.05%  "

# Row 50 - Algorand
$ws.Range("E50").Value = "  -2.32%  "

# Row 51 - EnergySwap
Set-TextValue $ws.Range("D51") "7.75"
$ws.Range("E51").Value = "  -0.08%  "
